##
## tests/utils/xlOil_Utils.xlsx
##
## Adds an "xloIndex" test section to the xloRef test sheet, renames the
## sheet xloRef -> xloRef-Index, and makes that sheet the active one.
##

$wb = $excel.ActiveWorkbook

$wsConcat = $wb.Worksheets.Item(2)   # "Concat-Split" - loses the selected/active sheet status
$wsRef    = $wb.Worksheets.Item(4)   # "xloRef" -> "xloRef-Index"

# ---------------------------------------------------------------------
# Rename the sheet and wipe its previous (small) test content so we can
# rebuild it to match the new, larger layout.
# ---------------------------------------------------------------------
$wsRef.Name = "xloRef-Index"
$wsRef.Cells.ClearContents()
$wsRef.Cells.ClearFormats()

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$wsRef.Columns.Item(2).ColumnWidth = 23.944010416666668

# ---------------------------------------------------------------------
# Header in K2 ("Data" shared-string, bold style copied from an existing
# bold header cell elsewhere in the workbook)
# ---------------------------------------------------------------------
$wb.Worksheets.Item(1).Range("J4").Copy()
$wsRef.Range("K2").PasteSpecial(-4122)
$wsRef.Range("K2").Value = "Data"

# ---------------------------------------------------------------------
# Section 1 - "xloRef: round trip"
# ---------------------------------------------------------------------
$wsRef.Range("B3").Value = "xloRef: round trip"

$wsRef.Range("K3").Value = 1
$wsRef.Range("L3").Value = 2
$wsRef.Range("M3").Value = 3
$wsRef.Range("K4").Value = 4
$wsRef.Range("L4").Value = 5
$wsRef.Range("M4").Value = 6
$wsRef.Range("K5").Value = 7
$wsRef.Range("L5").Value = 8
$wsRef.Range("M5").Value = 9

$wb.Worksheets.Item(2).Range("I5").Copy()
$wsRef.Range("K3:M5").PasteSpecial(-4122)

$wsRef.Range("E3").Formula = '=_xll.xloRef(K3:M5)'

$wb.Worksheets.Item(2).Range("I5").Copy()
$wsRef.Range("F3:H5").PasteSpecial(-4122)
$wsRef.Range("F3:H5").FormulaArray = '=_xll.xloVal(E3)'

$wsRef.Range("G3").Value = 2
$wsRef.Range("H3").Value = 3
$wsRef.Range("F4").Value = 4
$wsRef.Range("G4").Value = 5
$wsRef.Range("H4").Value = 6
$wsRef.Range("F5").Value = 7
$wsRef.Range("G5").Value = 8
$wsRef.Range("H5").Value = 9

$wsRef.Range("C3").FormulaArray = '=SUM(0+(F3:H5=K3:M5))=9'

# ---------------------------------------------------------------------
# Section 2 - "xloIndex: match Excel INDEX"
# ---------------------------------------------------------------------
$wsRef.Range("B8").Value = "xloIndex: match Excel INDEX"

$wsRef.Range("D8").Formula = '=INDEX($K$3:$M$5,2,2)'
$wsRef.Range("E8").Formula = '=_xll.xloIndex($E$3,2,2)'
$wsRef.Range("C8").Formula = '=D8=E8'

$wsRef.Range("D9").Formula = '=INDEX($K$3:$M$5,1,2)'
$wsRef.Range("E9").Formula = '=_xll.xloIndex($E$3,1,2)'
$wsRef.Range("C9").Formula = '=D9=E9'

$wsRef.Range("D10").Formula = '=INDEX(K3:K5,2)'
$wsRef.Range("E10").Formula = '=_xll.xloIndex(K3:K5,2)'
$wsRef.Range("C10").Formula = '=D10=E10'

# ---------------------------------------------------------------------
# Section 3 - "xloIndex: negative indices"
# ---------------------------------------------------------------------
$wsRef.Range("B13").Value = "xloIndex: negative indices"

$wsRef.Range("D13").Formula = '=INDEX(K3:M5,ROWS(K3:M5),COLUMNS(K3:M5))'
$wsRef.Range("E13").Formula = '=_xll.xloIndex($E$3,-1,-1)'
$wsRef.Range("C13").Formula = '=D13=E13'

# ---------------------------------------------------------------------
# Section 4 - "xloIndex: missing args"
# ---------------------------------------------------------------------
$wsRef.Range("B15").Value = "xloIndex: missing args"

$wsRef.Range("D15").Formula = '=K3'
$wsRef.Range("E15").Formula = '=L3'

$wb.Worksheets.Item(2).Range("I5").Copy()
$wsRef.Range("F15:G15").PasteSpecial(-4122)
$wsRef.Range("F15:G15").FormulaArray = '=_xll.xloIndex(E3,,,2,3)'
$wsRef.Range("G15").Value = 2

$wsRef.Range("C15").FormulaArray = '=SUM(0+(D15:E15=F15:G15))=2'

# ---------------------------------------------------------------------
# Section 5 - "xloIndex: zero toCol/toRow"
# ---------------------------------------------------------------------
$wsRef.Range("B17").Value = "xloIndex: zero toCol/toRow"

$wsRef.Range("D17").Formula = '=G4'
$wsRef.Range("E17").Formula = '=H4'

$wb.Worksheets.Item(2).Range("I5").Copy()
$wsRef.Range("F17:G18").PasteSpecial(-4122)
$wsRef.Range("F17:G18").FormulaArray = '=_xll.xloIndex($E$3,-2,-2, 0, 0)'
$wsRef.Range("G17").Value = 6
$wsRef.Range("F18").Value = 8
$wsRef.Range("G18").Value = 9

$wsRef.Range("D18").Formula = '=G5'
$wsRef.Range("E18").Formula = '=H5'

$wsRef.Range("C17").FormulaArray = '=SUM(0+(D17:E18=F17:G18))=4'

# ---------------------------------------------------------------------
# Page setup (print area / orientation) for the rebuilt sheet
# ---------------------------------------------------------------------
$wsRef.PageSetup.PaperSize = 9
$wsRef.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Sheet/selection bookkeeping: xloRef-Index becomes the active sheet and
# Concat-Split's previous selection/active status is cleared.
# ---------------------------------------------------------------------
$wsConcat.Range("L9").Select()

$wsRef.Activate()
$wsRef.Range("E23").Select()
